# Auto-generated edit script: applies scheduled-runner market-data refresh
# to the Asura_Profits price/profit columns (H-N) across all 8 sheets.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(135, 8).Value = 646.1
$ws.Cells.Item(135, 9).Value = 618.3333
$ws.Cells.Item(135, 11).Value = 5564.9997
$ws.Cells.Item(135, 13).Value = -3029.9997
$ws.Cells.Item(137, 8).Value = 2263.5366
$ws.Cells.Item(137, 9).Value = 1647.579
$ws.Cells.Item(137, 10).Value = 2795.5
$ws.Cells.Item(137, 11).Value = 4942.737
$ws.Cells.Item(137, 12).Value = 8386.5
$ws.Cells.Item(137, 13).Value = -2392.737
$ws.Cells.Item(137, 14).Value = -13486.5
$ws.Cells.Item(138, 8).Value = 3863.8767
$ws.Cells.Item(138, 9).Value = 3037.8235
$ws.Cells.Item(138, 10).Value = 4114.643
$ws.Cells.Item(138, 11).Value = 9113.470499999999
$ws.Cells.Item(138, 12).Value = 12343.929
$ws.Cells.Item(138, 13).Value = -3973.470499999999
$ws.Cells.Item(138, 14).Value = -22623.929
$ws.Cells.Item(141, 8).Value = 5218.2856
$ws.Cells.Item(141, 9).Value = 1707.8518
$ws.Cells.Item(141, 11).Value = 5123.555399999999
$ws.Cells.Item(141, 13).Value = 56.44460000000072

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 9018.612999999999
$ws.Cells.Item(32, 9).Value = 9586.493
$ws.Cells.Item(32, 10).Value = 5742.385
$ws.Cells.Item(32, 11).Value = 9586.493
$ws.Cells.Item(32, 12).Value = 5742.385
$ws.Cells.Item(32, 13).Value = -9299.493
$ws.Cells.Item(32, 14).Value = -6316.385
$ws.Cells.Item(61, 8).Value = 2586.5386
$ws.Cells.Item(61, 9).Value = 2460.5264
$ws.Cells.Item(61, 11).Value = 2460.5264
$ws.Cells.Item(61, 13).Value = -2248.5264
$ws.Cells.Item(74, 8).Value = 1753.7407
$ws.Cells.Item(74, 9).Value = 1643.0526
$ws.Cells.Item(74, 11).Value = 1643.0526
$ws.Cells.Item(74, 13).Value = -769.0526
$ws.Cells.Item(77, 8).Value = 1753.7407
$ws.Cells.Item(77, 9).Value = 1643.0526
$ws.Cells.Item(77, 11).Value = 8215.262999999999
$ws.Cells.Item(77, 13).Value = -3847.262999999999
$ws.Cells.Item(132, 8).Value = 8829.5625
$ws.Cells.Item(132, 9).Value = 8326.736999999999
$ws.Cells.Item(132, 10).Value = 9564.462
$ws.Cells.Item(132, 11).Value = 24980.211
$ws.Cells.Item(132, 12).Value = 28693.386
$ws.Cells.Item(132, 13).Value = -22450.211
$ws.Cells.Item(132, 14).Value = -33753.386
$ws.Cells.Item(136, 8).Value = 2586.5386
$ws.Cells.Item(136, 9).Value = 2460.5264
$ws.Cells.Item(136, 11).Value = 7381.5792
$ws.Cells.Item(136, 13).Value = -4831.5792

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(80, 8).Value = 1482034.2
$ws.Cells.Item(80, 9).Value = 3368020.2
$ws.Cells.Item(80, 10).Value = 188.14285
$ws.Cells.Item(80, 11).Value = 3368020.2
$ws.Cells.Item(80, 12).Value = 188.14285
$ws.Cells.Item(80, 13).Value = -3367022.2
$ws.Cells.Item(80, 14).Value = -2184.14285
$ws.Cells.Item(83, 8).Value = 1482034.2
$ws.Cells.Item(83, 9).Value = 3368020.2
$ws.Cells.Item(83, 10).Value = 188.14285
$ws.Cells.Item(83, 11).Value = 16840101
$ws.Cells.Item(83, 12).Value = 940.71425
$ws.Cells.Item(83, 13).Value = -16835109
$ws.Cells.Item(83, 14).Value = -10924.71425
$ws.Cells.Item(88, 8).Value = 28666.334
$ws.Cells.Item(88, 10).Value = 28666.334
$ws.Cells.Item(88, 12).Value = 28666.334
$ws.Cells.Item(88, 14).Value = -29478.334
$ws.Cells.Item(91, 8).Value = 28666.334
$ws.Cells.Item(91, 10).Value = 28666.334
$ws.Cells.Item(91, 12).Value = 28666.334
$ws.Cells.Item(91, 14).Value = -31474.334
$ws.Cells.Item(134, 8).Value = 2430
$ws.Cells.Item(134, 9).Value = 2172.3635
$ws.Cells.Item(134, 10).Value = 2902.3333
$ws.Cells.Item(134, 11).Value = 6517.0905
$ws.Cells.Item(134, 12).Value = 8706.999899999999
$ws.Cells.Item(134, 13).Value = -3982.0905
$ws.Cells.Item(134, 14).Value = -13776.9999

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2729.8538
$ws.Cells.Item(31, 9).Value = 2899.4
$ws.Cells.Item(31, 10).Value = 2568.3809
$ws.Cells.Item(31, 11).Value = 2899.4
$ws.Cells.Item(31, 12).Value = 2568.3809
$ws.Cells.Item(31, 13).Value = -2604.4
$ws.Cells.Item(31, 14).Value = -3158.3809
$ws.Cells.Item(34, 8).Value = 2729.8538
$ws.Cells.Item(34, 9).Value = 2899.4
$ws.Cells.Item(34, 10).Value = 2568.3809
$ws.Cells.Item(34, 11).Value = 2899.4
$ws.Cells.Item(34, 12).Value = 2568.3809
$ws.Cells.Item(34, 13).Value = -2697.4
$ws.Cells.Item(34, 14).Value = -2972.3809
$ws.Cells.Item(93, 8).Value = 4373
$ws.Cells.Item(93, 9).Value = 2710.3
$ws.Cells.Item(93, 10).Value = 21000
$ws.Cells.Item(93, 11).Value = 2710.3
$ws.Cells.Item(93, 12).Value = 21000
$ws.Cells.Item(93, 13).Value = -838.3000000000002
$ws.Cells.Item(93, 14).Value = -24744
$ws.Cells.Item(105, 8).Value = 3000
$ws.Cells.Item(105, 9).Value = 3000
$ws.Cells.Item(105, 10).Value = 0
$ws.Cells.Item(105, 11).Value = 3000
$ws.Cells.Item(105, 12).Value = 0
$ws.Cells.Item(105, 13).Value = -1253
$ws.Cells.Item(105, 14).ClearContents()

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(68, 8).Value = 182645.62
$ws.Cells.Item(68, 9).Value = 213398.48
$ws.Cells.Item(68, 10).Value = 1972.5
$ws.Cells.Item(68, 11).Value = 640195.4400000001
$ws.Cells.Item(68, 12).Value = 5917.5
$ws.Cells.Item(68, 13).Value = -639384.4400000001
$ws.Cells.Item(68, 14).Value = -7539.5
$ws.Cells.Item(71, 8).Value = 182645.62
$ws.Cells.Item(71, 9).Value = 213398.48
$ws.Cells.Item(71, 10).Value = 1972.5
$ws.Cells.Item(71, 11).Value = 1920586.32
$ws.Cells.Item(71, 12).Value = 17752.5
$ws.Cells.Item(71, 13).Value = -1916530.32
$ws.Cells.Item(71, 14).Value = -25864.5
$ws.Cells.Item(131, 8).Value = 23259218
$ws.Cells.Item(131, 9).Value = 561.36365
$ws.Cells.Item(131, 10).Value = 31254380
$ws.Cells.Item(131, 11).Value = 1684.09095
$ws.Cells.Item(131, 12).Value = 93763140
$ws.Cells.Item(131, 13).Value = 3355.90905
$ws.Cells.Item(131, 14).Value = -93773220
$ws.Cells.Item(132, 8).Value = 1346.0454
$ws.Cells.Item(132, 9).Value = 887.25
$ws.Cells.Item(132, 10).Value = 1608.2142
$ws.Cells.Item(132, 11).Value = 7985.25
$ws.Cells.Item(132, 12).Value = 14473.9278
$ws.Cells.Item(132, 13).Value = -5455.25
$ws.Cells.Item(132, 14).Value = -19533.9278

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(109, 8).Value = 11644.286
$ws.Cells.Item(109, 10).Value = 11644.286
$ws.Cells.Item(109, 12).Value = 11644.286
$ws.Cells.Item(109, 14).Value = -13724.286
$ws.Cells.Item(132, 8).Value = 3172.8462
$ws.Cells.Item(132, 9).Value = 2778
$ws.Cells.Item(132, 10).Value = 3804.6
$ws.Cells.Item(132, 11).Value = 8334
$ws.Cells.Item(132, 12).Value = 11413.8
$ws.Cells.Item(132, 13).Value = -5804
$ws.Cells.Item(132, 14).Value = -16473.8
$ws.Cells.Item(141, 8).Value = 55781.5
$ws.Cells.Item(141, 10).Value = 55781.5
$ws.Cells.Item(141, 12).Value = 55781.5
$ws.Cells.Item(141, 14).Value = -66141.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 3062.6365
$ws.Cells.Item(16, 9).Value = 746.8
$ws.Cells.Item(16, 10).Value = 4992.5
$ws.Cells.Item(16, 11).Value = 746.8
$ws.Cells.Item(16, 12).Value = 4992.5
$ws.Cells.Item(16, 13).Value = -576.8
$ws.Cells.Item(16, 14).Value = -5332.5
$ws.Cells.Item(46, 8).Value = 902.8570999999999
$ws.Cells.Item(46, 10).Value = 970
$ws.Cells.Item(46, 12).Value = 970
$ws.Cells.Item(46, 14).Value = -1346
$ws.Cells.Item(96, 8).Value = 54500
$ws.Cells.Item(96, 10).Value = 54500
$ws.Cells.Item(96, 12).Value = 54500
$ws.Cells.Item(96, 14).Value = -59992
$ws.Cells.Item(132, 8).Value = 4180.9067
$ws.Cells.Item(132, 9).Value = 4210.1353
$ws.Cells.Item(132, 10).Value = 4000.6667
$ws.Cells.Item(132, 11).Value = 12630.4059
$ws.Cells.Item(132, 12).Value = 12002.0001
$ws.Cells.Item(132, 13).Value = -10100.4059
$ws.Cells.Item(132, 14).Value = -17062.0001
$ws.Cells.Item(141, 8).Value = 50000
$ws.Cells.Item(141, 10).Value = 50000
$ws.Cells.Item(141, 12).Value = 50000
$ws.Cells.Item(141, 14).Value = -60360

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(109, 8).Value = 24895.5
$ws.Cells.Item(109, 10).Value = 24895.5
$ws.Cells.Item(109, 12).Value = 24895.5
$ws.Cells.Item(109, 14).Value = -27669.5
$ws.Cells.Item(132, 8).Value = 2299.0256
$ws.Cells.Item(132, 9).Value = 2380.7144
$ws.Cells.Item(132, 10).Value = 2203.7222
$ws.Cells.Item(132, 11).Value = 7142.1432
$ws.Cells.Item(132, 12).Value = 6611.1666
$ws.Cells.Item(132, 13).Value = -4612.1432
$ws.Cells.Item(132, 14).Value = -11671.1666
